# Commit: "10 part of method is ended"
#
# The component names in the first column of the binary-interaction
# coefficient matrix (rows 2-4) are updated: CH4 -> CH41, C2H6 -> C2H61,
# C3H8 -> C3H81. These are brand-new label strings (added to the shared
# string table), distinct from the existing "CH4"/"C2H6"/"C3H8" headers
# still used elsewhere (e.g. row 1 / column headers), so the header row
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CH41"
$ws.Range("A3").Value = "C2H61"
$ws.Range("A4").Value = "C3H81"
